$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.024.77'
$ws.Range('E2').Value = '  -2.18%  '
$ws.Range('D3').Value = '1.667.29'
$ws.Range('E3').Value = '  -1.53%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.89'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.32%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5110'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.23%  '
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2656'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06402'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.80'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07438'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.11%  '
$ws.Range('D12').Value = '1.668.37'
$ws.Range('E12').Value = '  -1.60%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.512'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.17%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5841'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.92%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008584'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.74%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.29'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.66%  '
$ws.Range('D17').Value = '26.071.17'
$ws.Range('E17').Value = '  -2.10%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.937'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.04%  '
$ws.Range('E19').Value = '  -0.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.77'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.70%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '191.97'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.208'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.83%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.005'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '144.71'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '7.619'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.78%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1198'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.70%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.67'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.92%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06524'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +13.63%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.340'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.318'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.83%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.539'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.48%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.516'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.22%  '
$ws.Range('E33').Value = '  +0.61%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.019'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6122'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.02%  '
$ws.Range('E36').Value = '  +0.29%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.683'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.276'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +7.62%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01603'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.07%  '
$ws.Range('D40').Value = '1.089.73'
$ws.Range('E40').Value = '  -0.91%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8642'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.60%  '
$ws.Range('E42').Value = '  +0.22%  '
$ws.Range('E43').Value = '  +1.47%  '
$ws.Range('D44').Value = '1.816.44'
$ws.Range('E44').Value = '  -1.92%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000114'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.40%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.40'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.30%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.007'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.16%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.050'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.55%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05230'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.09%  '
$ws.Range('B50').Value = 'Aptos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.086'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.07%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4285'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.95%  '
